$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row: insert a new "Residual deviance" column before the AIC column ---
# Copy C1 ("AIC", with its header style/border) into D1, then retext C1.
$ws.Range("C1").Copy($ws.Range("D1"))
$ws.Range("C1").Value = "Residual deviance"

# --- Update Segment 1 / Segment 2 labels in column A ---
$ws.Range("A4").Value = "Segment 1: <0.0.3698 mg P/L"

$ws.Range("A5").Value = "Segment 2: X0.03698 mg P/L"
$seg2ge = $ws.Range("A5").Characters(12, 1)
$seg2ge.Font.Name = "Arial"
$seg2ge.Font.Size = 12
$seg2rest = $ws.Range("A5").Characters(13, 14)
$seg2rest.Font.Name = "Times New Roman"
$seg2rest.Font.Size = 12

# --- Numeric data: Residual deviance (C) and AIC (D) for every model row ---
$ws.Range("C2").Value = 43.713
$ws.Range("D2").Value = 69.398

$ws.Range("C3").Value = 42.384
$ws.Range("D3").Value = 69.091
$ws.Range("E3").Value = "NOT ADDITIVE "

$ws.Range("C4").Value = 8.773
$ws.Range("D4").Value = 20.846

$ws.Range("C5").Value = 31.318
$ws.Range("D5").Value = 46.261

$ws.Range("C6").Formula = "=SUM(C7:C8)"
$ws.Range("D6").Formula = "=SUM(D7:D8)"
$ws.Range("E6").Value = "UNSURE IF THESE SHOULD BE ADDED "

$ws.Range("C7").Value = 3.22
$ws.Range("D7").Value = 6.4106

$ws.Range("C8").Value = 8.5632
$ws.Range("D8").Value = 20.679

# --- Number format + centering for the new numeric columns ---
$ws.Range("C2:C5").HorizontalAlignment = -4108
$ws.Range("C8").HorizontalAlignment = -4108
$ws.Range("D2:D5").HorizontalAlignment = -4108
$ws.Range("D7:D8").HorizontalAlignment = -4108

$ws.Range("C6").NumberFormat = "0.000"
$ws.Range("C6").HorizontalAlignment = -4108
$ws.Range("D6").NumberFormat = "0.000"
$ws.Range("D6").HorizontalAlignment = -4108
$ws.Range("C7").NumberFormat = "0.000"
$ws.Range("C7").HorizontalAlignment = -4108

# --- Column C width ---
$ws.Columns("C").ColumnWidth = 10.25

# --- Selection on Sheet1 ---
$ws.Range("C11:C12").Select()

# --- New "Appendix E" worksheet, placed after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "Appendix E"
$ws2.Range("E36").Select()

# Re-activate Sheet1 so it remains the selected tab on save.
$ws.Activate()
